$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "26.602.94"  # D2: '26.605.82' -> '26.602.94'
Set-TextValue $ws.Cells.Item(2, 5) "  +0.44%  "  # E2: '  +0.46%  ' -> '  +0.44%  '

Set-TextValue $ws.Cells.Item(3, 4) "1.738.94"  # D3: '1.739.03' -> '1.738.94'
Set-TextValue $ws.Cells.Item(3, 5) "  +0.68%  "  # E3: '  +0.66%  ' -> '  +0.68%  '

Set-TextValue $ws.Cells.Item(4, 5) "  +0.03%  "  # E4: '  +0.04%  ' -> '  +0.03%  '

Set-TextValue $ws.Cells.Item(5, 4) "246.07"  # D5: '246.06' -> '246.07'
Set-TextValue $ws.Cells.Item(5, 5) "  +0.23%  "  # E5: '  +0.06%  ' -> '  +0.23%  '

Set-TextValue $ws.Cells.Item(6, 4) "1.0000"  # D6: '0.9999' -> '1.0000'

Set-TextValue $ws.Cells.Item(7, 4) "0.4954"  # D7: '0.4963' -> '0.4954'
Set-TextValue $ws.Cells.Item(7, 5) "  +3.26%  "  # E7: '  +3.49%  ' -> '  +3.26%  '

Set-TextValue $ws.Cells.Item(8, 4) "0.2675"  # D8: '0.2677' -> '0.2675'
Set-TextValue $ws.Cells.Item(8, 5) "  -0.34%  "  # E8: '  -0.35%  ' -> '  -0.34%  '

Set-TextValue $ws.Cells.Item(9, 4) "0.06265"  # D9: '0.06270' -> '0.06265'
Set-TextValue $ws.Cells.Item(9, 5) "  +0.64%  "  # E9: '  +0.67%  ' -> '  +0.64%  '

Set-TextValue $ws.Cells.Item(10, 4) "1.740.54"  # D10: '1.741.57' -> '1.740.54'
Set-TextValue $ws.Cells.Item(10, 5) "  +0.73%  "  # E10: '  +0.86%  ' -> '  +0.73%  '

Set-TextValue $ws.Cells.Item(11, 4) "0.07045"  # D11: '0.07042' -> '0.07045'

Set-TextValue $ws.Cells.Item(12, 5) "  +0.26%  "  # E12: '  +0.07%  ' -> '  +0.26%  '

Set-TextValue $ws.Cells.Item(13, 2) "Polygon"  # B13: 'Polkadot' -> 'Polygon'
Set-TextValue $ws.Cells.Item(13, 3) "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"  # C13: 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' -> 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Cells.Item(13, 4) "0.6129"  # D13: '4.588' -> '0.6129'
Set-TextValue $ws.Cells.Item(13, 5) "  -1.03%  "  # E13: '  +1.51%  ' -> '  -1.03%  '

Set-TextValue $ws.Cells.Item(14, 2) "Polkadot"  # B14: 'Polygon' -> 'Polkadot'
Set-TextValue $ws.Cells.Item(14, 3) "https://coinranking.com/coin/25W7FG7om+polkadot-dot"  # C14: 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' -> 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Cells.Item(14, 4) "4.586"  # D14: '0.6130' -> '4.586'
Set-TextValue $ws.Cells.Item(14, 5) "  +1.55%  "  # E14: '  -1.22%  ' -> '  +1.55%  '

Set-TextValue $ws.Cells.Item(15, 4) "77.95"  # D15: '77.98' -> '77.95'
Set-TextValue $ws.Cells.Item(15, 5) "  +1.06%  "  # E15: '  +1.02%  ' -> '  +1.06%  '

Set-TextValue $ws.Cells.Item(17, 4) "26.619.30"  # D17: '26.622.62' -> '26.619.30'
Set-TextValue $ws.Cells.Item(17, 5) "  +0.45%  "  # E17: '  +0.47%  ' -> '  +0.45%  '

Set-TextValue $ws.Cells.Item(18, 5) "  +0.01%  "  # E18: '  +0.05%  ' -> '  +0.01%  '

Set-TextValue $ws.Cells.Item(19, 4) "0.000007264"  # D19: '0.000007246' -> '0.000007264'
Set-TextValue $ws.Cells.Item(19, 5) "  +4.68%  "  # E19: '  +4.25%  ' -> '  +4.68%  '

Set-TextValue $ws.Cells.Item(20, 4) "11.57"  # D20: '11.55' -> '11.57'
Set-TextValue $ws.Cells.Item(20, 5) "  -0.91%  "  # E20: '  -1.04%  ' -> '  -0.91%  '

Set-TextValue $ws.Cells.Item(21, 4) "1.968.16"  # D21: '1.969.10' -> '1.968.16'
Set-TextValue $ws.Cells.Item(21, 5) "  +0.91%  "  # E21: '  +0.96%  ' -> '  +0.91%  '

Set-TextValue $ws.Cells.Item(22, 4) "4.557"  # D22: '4.562' -> '4.557'
Set-TextValue $ws.Cells.Item(22, 5) "  +0.46%  "  # E22: '  +0.56%  ' -> '  +0.46%  '

Set-TextValue $ws.Cells.Item(23, 4) "8.720"  # D23: '8.717' -> '8.720'
Set-TextValue $ws.Cells.Item(23, 5) "  -2.68%  "  # E23: '  -2.75%  ' -> '  -2.68%  '

Set-TextValue $ws.Cells.Item(24, 4) "5.267"  # D24: '5.272' -> '5.267'
Set-TextValue $ws.Cells.Item(24, 5) "  -0.54%  "  # E24: '  -0.47%  ' -> '  -0.54%  '

Set-TextValue $ws.Cells.Item(25, 4) "138.83"  # D25: '138.85' -> '138.83'
Set-TextValue $ws.Cells.Item(25, 5) "  +1.78%  "  # E25: '  +1.83%  ' -> '  +1.78%  '

Set-TextValue $ws.Cells.Item(26, 5) "  +0.35%  "  # E26: '  +0.28%  ' -> '  +0.35%  '

Set-TextValue $ws.Cells.Item(27, 4) "1.424"  # D27: '1.422' -> '1.424'
Set-TextValue $ws.Cells.Item(27, 5) "  +1.12%  "  # E27: '  +1.05%  ' -> '  +1.12%  '

Set-TextValue $ws.Cells.Item(28, 4) "1.762"  # D28: '1.759' -> '1.762'
Set-TextValue $ws.Cells.Item(28, 5) "  -2.27%  "  # E28: '  -2.42%  ' -> '  -2.27%  '

Set-TextValue $ws.Cells.Item(29, 4) "107.22"  # D29: '107.27' -> '107.22'
Set-TextValue $ws.Cells.Item(29, 5) "  +0.44%  "  # E29: '  +0.48%  ' -> '  +0.44%  '

Set-TextValue $ws.Cells.Item(30, 5) "  +1.31%  "  # E30: '  +1.26%  ' -> '  +1.31%  '

Set-TextValue $ws.Cells.Item(31, 4) "0.08048"  # D31: '0.08049' -> '0.08048'
Set-TextValue $ws.Cells.Item(31, 5) "  +0.32%  "  # E31: '  +0.28%  ' -> '  +0.32%  '

Set-TextValue $ws.Cells.Item(32, 4) "3.732"  # D32: '3.730' -> '3.732'
Set-TextValue $ws.Cells.Item(32, 5) "  +0.23%  "  # E32: '  +0.08%  ' -> '  +0.23%  '

Set-TextValue $ws.Cells.Item(33, 4) "0.04622"  # D33: '0.04621' -> '0.04622'
Set-TextValue $ws.Cells.Item(33, 5) "  +1.18%  "  # E33: '  +1.24%  ' -> '  +1.18%  '

Set-TextValue $ws.Cells.Item(35, 4) "2.613"  # D35: '2.612' -> '2.613'
Set-TextValue $ws.Cells.Item(35, 5) "  -0.17%  "  # E35: '  -0.19%  ' -> '  -0.17%  '

Set-TextValue $ws.Cells.Item(36, 4) "1.014"  # D36: '1.013' -> '1.014'
Set-TextValue $ws.Cells.Item(36, 5) "  +2.55%  "  # E36: '  +2.22%  ' -> '  +2.55%  '

Set-TextValue $ws.Cells.Item(37, 4) "0.6373"  # D37: '0.6395' -> '0.6373'
Set-TextValue $ws.Cells.Item(37, 5) "  +0.10%  "  # E37: '  +0.31%  ' -> '  +0.10%  '

Set-TextValue $ws.Cells.Item(38, 4) "2.061"  # D38: '2.067' -> '2.061'
Set-TextValue $ws.Cells.Item(38, 5) "  -2.39%  "  # E38: '  -1.72%  ' -> '  -2.39%  '

Set-TextValue $ws.Cells.Item(39, 4) "0.9020"  # D39: '0.9010' -> '0.9020'
Set-TextValue $ws.Cells.Item(39, 5) "  -3.73%  "  # E39: '  -3.82%  ' -> '  -3.73%  '

Set-TextValue $ws.Cells.Item(40, 4) "2.423"  # D40: '2.425' -> '2.423'
Set-TextValue $ws.Cells.Item(40, 5) "  +0.46%  "  # E40: '  +0.44%  ' -> '  +0.46%  '

Set-TextValue $ws.Cells.Item(41, 5) "  -0.01%  "  # E41: '  -0.34%  ' -> '  -0.01%  '

Set-TextValue $ws.Cells.Item(43, 4) "101.93"  # D43: '101.95' -> '101.93'
Set-TextValue $ws.Cells.Item(43, 5) "  -2.88%  "  # E43: '  -2.82%  ' -> '  -2.88%  '

Set-TextValue $ws.Cells.Item(44, 4) "5.438"  # D44: '5.439' -> '5.438'
Set-TextValue $ws.Cells.Item(44, 5) "  -4.57%  "  # E44: '  -4.62%  ' -> '  -4.57%  '

Set-TextValue $ws.Cells.Item(45, 4) "0.3934"  # D45: '0.3935' -> '0.3934'
Set-TextValue $ws.Cells.Item(45, 5) "  +0.94%  "  # E45: '  +0.61%  ' -> '  +0.94%  '

Set-TextValue $ws.Cells.Item(46, 4) "6.857"  # D46: '6.859' -> '6.857'
Set-TextValue $ws.Cells.Item(46, 5) "  -1.22%  "  # E46: '  -1.54%  ' -> '  -1.22%  '

Set-TextValue $ws.Cells.Item(47, 4) "0.1181"  # D47: '0.1182' -> '0.1181'
Set-TextValue $ws.Cells.Item(47, 5) "  -0.74%  "  # E47: '  -0.71%  ' -> '  -0.74%  '

Set-TextValue $ws.Cells.Item(48, 5) "  +1.10%  "  # E48: '  +1.14%  ' -> '  +1.10%  '

Set-TextValue $ws.Cells.Item(49, 4) "30.69"  # D49: '30.68' -> '30.69'
Set-TextValue $ws.Cells.Item(49, 5) "  -1.12%  "  # E49: '  -1.20%  ' -> '  -1.12%  '

Set-TextValue $ws.Cells.Item(50, 4) "7.805"  # D50: '7.823' -> '7.805'
Set-TextValue $ws.Cells.Item(50, 5) "  -0.86%  "  # E50: '  -1.30%  ' -> '  -0.86%  '

Set-TextValue $ws.Cells.Item(51, 4) "1.252"  # D51: '1.254' -> '1.252'
Set-TextValue $ws.Cells.Item(51, 5) "  -1.14%  "  # E51: '  -1.00%  ' -> '  -1.14%  '
